$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G ("K") values for rows 2, 3, 5, 6 to reflect the
# regenerated save_data (switch from Strike# to K, recalculated values).
$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 4
$ws.Range("G5").Value = 3
$ws.Range("G6").Value = 2
